$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column before column N (shifts old N:R -> O:S).
# Excel copies the left-neighbour's (M) formatting into the new column
# automatically, which already gives the correct header style for N1.
$ws.Range("N1").EntireColumn.Insert()

# New "Die Making" column header + values
$ws.Range("N1").Value = "Die Making"

$ws.Range("N2").Value = 1500
$ws.Range("N3").Value = 2000
$ws.Range("N4").Value = 3000
$ws.Range("N5").Value = 3500
$ws.Range("N6").Value = 4000
$ws.Range("N7").Value = 5000

# Match the formatting already used by the (shifted) price column O4:O7
$ws.Range("O4").Copy()
$ws.Range("N2:N7").PasteSpecial(-4122)

# Row 2 value corrections for the shifted coating columns
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 3
$ws.Range("R2").Value = 15

# Update selection to match the post-edit state
$ws.Range("N8").Select()
